$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply date format (yyyy-mm-dd, escaped dashes) to the existing date
#     column so new rows created below inherit matching formatting ---
$ws.Range("A2:A6").NumberFormat = "yyyy\-mm\-dd"

# --- Fill in "Beat Vegas?" results for the already-played games (Jan 5) ---
$ws.Range("G2").Value = "No"
$ws.Range("G3").Value = "No"
$ws.Range("G4").Value = "Yes"
$ws.Range("G5").Value = "Yes"
$ws.Range("G6").Value = "Yes"

# --- Copy the formatting (date number format, etc.) from the last
#     existing data row down across the new rows (7-17). Column G is
#     intentionally excluded since none of the new games have a
#     "Beat Vegas?" verdict yet. ---
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A7:F17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New games run through the model for Jan 6, 2021 (serial 44202) ---
$newRows = @(
  @(44202, "IND", "HOU", -2.5,  -3,    0.5),
  @(44202, "PHI", "WAS", -6.5,   0.1, -6.6),
  @(44202, "ORL", "CLE", -6,    -6.2,  0.20000000000000021),
  @(44202, "MIA", "BOS", -2.5,   0.8, -3.3),
  @(44202, "ATL", "CHO", -6,    -7.8,  1.8),
  @(44202, "NYK", "UTA",  7.5,   6.3,  1.2),
  @(44202, "NOP", "OKC", -8,     2.9, -10.9),
  @(44202, "MIL", "DET", -12,  -22.6, 10.6),
  @(44202, "PHO", "TOR", -3,   -10.4,  7.4),
  @(44202, "SAC", "CHI", -7,     7.2, -14.2),
  @(44202, "GSW", "LAC", -1,    -3.5,  2.5)
)

$r = 7
foreach ($row in $newRows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $r++
}

# --- Resize column A to fit the new date strings ---
$ws.Columns.Item(1).ColumnWidth = 9.5
